$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns per latest data pull.
# D-column prices are written with a leading apostrophe to force text storage
# (matching the sheet's original text cells) and then the cell style is reset
# to "Normal" so the quote-prefix flag does not linger on the cell.

$ws.Range("D2").Value = "'40.238.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.14%  "

$ws.Range("D3").Value = "'2.247.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("D5").Value = "'294.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").Value = "'87.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.69%  "

$ws.Range("D7").Value = "'0.516"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.94%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "

$ws.Range("D10").Value = "'31.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.26%  "

$ws.Range("D11").Value = "'0.0799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.71%  "

$ws.Range("D12").Value = "'47.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("E13").Value = "  +0.97%  "

$ws.Range("E14").Value = "  +6.33%  "

$ws.Range("D15").Value = "'2.597.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").Value = "'14.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").Value = "'2.254.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").Value = "'0.741"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.72%  "

$ws.Range("D19").Value = "'40.164.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.16%  "

$ws.Range("E20").Value = "  +4.26%  "

$ws.Range("D21").Value = "'5.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.05%  "

$ws.Range("D22").Value = "'10.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.19%  "

$ws.Range("D23").Value = "'65.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").Value = "'236.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.02%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  +3.77%  "

$ws.Range("E27").Value = "  +7.84%  "

$ws.Range("D28").Value = "'23.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.16%  "

$ws.Range("D29").Value = "'2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.30%  "

$ws.Range("D30").Value = "'9.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.66%  "

$ws.Range("D31").Value = "'33.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.93%  "

$ws.Range("D32").Value = "'153.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "'4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.07%  "

$ws.Range("D35").Value = "'0.0720"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.32%  "

$ws.Range("D36").Value = "'2.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "

$ws.Range("D37").Value = "'16.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.31%  "

$ws.Range("E38").Value = "  +6.99%  "

$ws.Range("E39").Value = "  +2.74%  "

$ws.Range("D40").Value = "'2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.82%  "

$ws.Range("E41").Value = "  +6.28%  "

$ws.Range("D42").Value = "'3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.31%  "

$ws.Range("D43").Value = "'2.023.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.24%  "

$ws.Range("E44").Value = "  +7.72%  "

$ws.Range("D45").Value = "'0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.57%  "

$ws.Range("D46").Value = "'10.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.63%  "

$ws.Range("D47").Value = "'16.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("D48").Value = "'2.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.36%  "

$ws.Range("D49").Value = "'2.478.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.61%  "

$ws.Range("D50").Value = "'71.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.76%  "

$ws.Range("D51").Value = "'1.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.15%  "
